$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regen save_data to use K (strikeouts) instead of the old Strike# values.
# Column G holds "K" per the header row; update the per-game values to the
# freshly computed/regenerated numbers.
$newK = @{
    2  = 1
    3  = 3
    4  = 4
    5  = 4
    6  = 1
    7  = 0
    8  = 6
    9  = 5
    10 = 1
    11 = 1
    12 = 6
    13 = 3
    15 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
